$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Initial_DataSetUp")

# Update "Is Parallel Execution" field value from Yes to No
$ws.Range("B5").Value = "No"

# Update "Parallel Execution Count" field value from 2 to 1 (kept as text,
# not a number, matching the original cell type/style)
$b6 = $ws.Range("B6")
$b6.Formula = '="1"'
$b6.Copy()
$b6.PasteSpecial(-4163)  # xlPasteValues

# Move the active selection to B6, matching the last-edited cell
$ws.Activate()
$ws.Range("B6").Select()
